$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H18").Value = "CPU"
$ws.Range("I18").Value = 30

$ws.Range("H19").Value = "GPG/CPU"
$ws.Range("I19").Formula = "=I18/I12"

$ws.Range("H20").Select()
